# Fix css class name and add test files.
# - Correct the header labels in the "variable_mapping" sheet (remove the
#   stray space that had crept into "Platform Name" / "Header Roow").
# - Normalize the formatting of the "Kakao" row (row 6) and column B so
#   they use the same font/number-format as the rest of the table instead
#   of the stray duplicate "theme color" font that only that row had.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variable_mapping")

# 1) Fix the header text typos.
$ws.Range("A1").Value = "PlatformName"
$ws.Range("B1").Value = "HeaderRoow"

# 2) Re-apply the same font used by the other data rows (2-5, 7) to row 6,
#    which previously carried its own slightly-different (duplicate) font.
#    Number formats already match (General for C:O, #,##0 for B) -- only
#    the font differs, so only touch the font here.
$sourceFont = $ws.Range("B2").Font
$targetRow = $ws.Range("B6:O6")
$targetRow.Font.Color = $sourceFont.Color
$targetRow.Font.Name = $sourceFont.Name
$targetRow.Font.Size = $sourceFont.Size

# 3) Column B as a whole should carry the same (now de-duplicated) font.
$ws.Columns.Item(2).Font.Color = $sourceFont.Color
$ws.Columns.Item(2).Font.Name = $sourceFont.Name
$ws.Columns.Item(2).Font.Size = $sourceFont.Size

# 4) Row 6's height changes from 18 to 19.5 (matching the other rows).
$ws.Rows.Item(6).RowHeight = 19.5
